$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDSheet")

# Reword five vehicle registration entries (reorder the plate/registration
# tokens) on the "Основное средство" column (B). Assignment order matters:
# it controls the order new entries are appended to the shared-string table.
$ws.Range("B42").Value = "BMW X3 #О987УС120У33, С234ОР799"
$ws.Range("B49").Value = "Mercedes-Benz GLE #С789ЕС120Е33, ddd00ef, О123ЕР799"
$ws.Range("B33").Value = "Toyota Prius #У797АА120С33, О777ОР799, #D797QQ120Е33"
$ws.Range("B23").Value = "Ford Mustang D123SS799, #С361СС920У33, Р123РХ799"
$ws.Range("B20").Value = "Chevrolet Silverado #Р678СО120Н22, Е456СР799"

# Update the saved view state: scroll position and active selection.
$ws.Range("F40").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
